# Insert a new slide between the existing slide 3 ("SMART Research
# Question") and slide 4 ("Reshaping Our Data for the SMART Q"), using the
# same "Title and Content" layout already used elsewhere in the deck.
# This pushes the old slide 4 down to become slide 5, and the new slide
# becomes the new slide 4 with just a title placeholder filled in.
$p = $ppt.ActivePresentation

$layouts = $p.SlideMaster.CustomLayouts
$titleAndContent = $null
for ($i = 1; $i -le $layouts.Count; $i++) {
    if ($layouts.Item($i).Name -eq "Title and Content") {
        $titleAndContent = $layouts.Item($i)
        break
    }
}

$newSlide = $p.Slides.AddSlide(4, $titleAndContent)

# Fill in the title placeholder; leave the content placeholder empty.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "RAJEEV START HERE"
